# Consolidated Transmitter Board BOM with master BOM
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Update Manufacturer Part Number 1 values (column C) for the rows whose
# part numbers were reconciled against the master BOM.
$ws.Range("C6").Value = "CL10B105MO8NNWC"
$ws.Range("C11").Value = "GRM21BR61E106KA73K"
$ws.Range("C24").Value = "RC0603JR-0710KL"
$ws.Range("C35").Value = "ESR10EZPJ681"

# Update Supplier Unit Price 1 (G) and Supplier Subtotal 1 (H) to match the
# pricing from the consolidated master BOM.
$ws.Range("G6").Value = 0.008
$ws.Range("H6").Value = 0.032

$ws.Range("G11").Value = 0.104
$ws.Range("H11").Value = 1.04

$ws.Range("G24").Value = 0.015
$ws.Range("H24").Value = 0.15

$ws.Range("G35").Value = 0.1
$ws.Range("H35").Value = 0.1
